$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 511
$ws1.Range("F4").Value = 474
$ws1.Range("F5").Value = 8856
$ws1.Range("F7").Value = 11353
$ws1.Range("G20").Value = 68
$ws1.Range("F22").Value = 1925
$ws1.Range("F23").Value = 732
$ws1.Range("F24").Value = 664
$ws1.Range("F25").Value = 367
$ws1.Range("F27").Value = 81
$ws1.Range("F28").Value = 615
$ws1.Range("F30").Value = 1358
$ws1.Range("F32").Value = 14
$ws1.Range("F33").Value = 8
$ws1.Range("F35").Value = 472
$ws1.Range("F36").Value = 314
$ws1.Range("F37").Value = 28
$ws1.Range("F39").Value = 338
$ws1.Range("F43").Value = 401
$ws1.Range("F45").Value = 822
$ws1.Range("F46").Value = 661
$ws1.Range("F48").Value = 179
$ws1.Range("F49").Value = 167

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F18").Value = 73
$ws2.Range("F19").Value = 115

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2845

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 511
$ws4.Range("F8").Value = 474
$ws4.Range("F9").Value = 8856
$ws4.Range("F11").Value = 11353
$ws4.Range("G19").Value = 68
$ws4.Range("F20").Value = 1925
$ws4.Range("F21").Value = 732
$ws4.Range("F22").Value = 664
$ws4.Range("F23").Value = 367
$ws4.Range("F26").Value = 615
$ws4.Range("F29").Value = 1358
$ws4.Range("F30").Value = 14
$ws4.Range("F35").Value = 472
$ws4.Range("F36").Value = 314
$ws4.Range("F37").Value = 73
$ws4.Range("F41").Value = 401
$ws4.Range("F43").Value = 822
$ws4.Range("F46").Value = 661
$ws4.Range("F48").Value = 179
$ws4.Range("F49").Value = 167
